# Add a new order line (row 4) to the completed order sheet, matching the
# existing sheet's convention of storing every value as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "'173339"
$ws.Range("B4").Value = "Chobani - Drinkable Yogurt"
$ws.Range("C4").Value = "'5"
$ws.Range("D4").Value = "'17.99"
$ws.Range("E4").Value = "'89.95"
